# Generate Report for Handback
# ------------------------------------------------------------------
# The localization files have now been handed back from the vendor:
#   * Status moves from "Ready for handoff" -> "Handed back: in sync with en-US"
#   * The per-language "Latest Target File" / "Latest Handback File" /
#     "Latest Handback DateTime" columns get populated for both rows
#     on both the zh-cn and de-de detail sheets.
#   * A handful of columns are widened so the new (longer) values fit.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/7923e81af238441569c8e21c80afdce11cc6fbdc/e2e/"
$mdA = "a2ae3aa4-d209-4a02-a640-e6c65e5ba926.md"
$mdF = "f6537433-6b8c-4357-a94d-84ed17d31411.md"

# ------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    This string is shared across Overview (E/F), zh-cn (C) and de-de (C),
#    so updating every cell that currently shows it keeps them all in sync.
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

# ------------------------------------------------------------------
# 2. zh-cn detail sheet: fill in Latest Target File (I), Latest Handback
#    File (J) and stamp Latest Handback DateTime (K).
# ------------------------------------------------------------------
$wsZh.Range("I2").Value = $mdA
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($baseUrl + $mdA), "", "", $mdA)
$wsZh.Range("J2").Value = "a2ae3aa4-d209-4a02-a640-e6c65e5ba926.78e5269e4ef310fca9491d991caeef7b2ee1626f.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-09-03 08:32:51"

$wsZh.Range("I3").Value = $mdF
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($baseUrl + $mdF), "", "", $mdF)
$wsZh.Range("J3").Value = "f6537433-6b8c-4357-a94d-84ed17d31411.821eedb5654191eaae5ec64288ecbdc3c8609a73.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-09-03 08:32:51"

# ------------------------------------------------------------------
# 3. de-de detail sheet: same shape, different target files / timestamp.
# ------------------------------------------------------------------
$wsDe.Range("I2").Value = $mdA
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($baseUrl + $mdA), "", "", $mdA)
$wsDe.Range("J2").Value = "a2ae3aa4-d209-4a02-a640-e6c65e5ba926.78e5269e4ef310fca9491d991caeef7b2ee1626f.de-de.xlf"
$wsDe.Range("K2").Value = "2016-09-03 08:32:57"

$wsDe.Range("I3").Value = $mdF
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($baseUrl + $mdF), "", "", $mdF)
$wsDe.Range("J3").Value = "f6537433-6b8c-4357-a94d-84ed17d31411.821eedb5654191eaae5ec64288ecbdc3c8609a73.de-de.xlf"
$wsDe.Range("K3").Value = "2016-09-03 08:32:57"

# ------------------------------------------------------------------
# 4. Widen columns to fit the newly-populated / longer text.
#    ColumnWidth is in "characters"; Excel stores the underlying sheet
#    width in pixel (1/6 character) steps, so we feed it values that
#    land the stored width on the target pixel boundary.
# ------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668   # -> stored width 30 (was ~17.2)
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668   # -> stored width 30 (was ~17.2)

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668         # Status column -> stored width 30
$wsZh.Columns.Item(9).ColumnWidth = 39.16666666666667          # Latest Target File -> stored width 40
$wsZh.Columns.Item(10).ColumnWidth = 39.16666666666667         # Latest Handback File -> stored width 40

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668         # Status column -> stored width 30
$wsDe.Columns.Item(9).ColumnWidth = 39.16666666666667          # Latest Target File -> stored width 40
$wsDe.Columns.Item(10).ColumnWidth = 39.16666666666667         # Latest Handback File -> stored width 40
